$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for every existing data row
# (rows 2-339) from 2023-10-03 (45202) to 2023-10-04 (45203).
$ws.Range("C2:C339").Value = 45203

# Row 339 gains an explicit row-height / customHeight flag in the target file.
$ws.Rows.Item(339).RowHeight = 15

# Append the new record as row 340.
$ws.Range("A340").Value = "A 46997-2023"

$ws.Range("B340").Value = 45201
$ws.Range("B340").NumberFormat = "YYYY-MM-DD"

$ws.Range("C340").Value = 45203
$ws.Range("C340").NumberFormat = "YYYY-MM-DD"

$ws.Range("D340").Value = "ÖREBRO LÄN"
$ws.Range("E340").Value = "HÄLLEFORS"
$ws.Range("F340").Value = "Bergvik skog väst AB"

$ws.Range("G340").Value = 0.9
$ws.Range("H340").Value = 0
$ws.Range("I340").Value = 0
$ws.Range("J340").Value = 0
$ws.Range("K340").Value = 0
$ws.Range("L340").Value = 0
$ws.Range("M340").Value = 0
$ws.Range("N340").Value = 0
$ws.Range("O340").Value = 0
$ws.Range("P340").Value = 0
$ws.Range("Q340").Value = 0

$ws.Range("R340").Value = ""
$ws.Range("R340").WrapText = $true

Write-Output "edit applied"
